$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4146968223820794
$ws.Range("C2").Value = 0.9918227464331973
$ws.Range("D2").Value = 0.4952318647003984
$ws.Range("F2").Value = "Pipeline(steps=[('model', RandomForestRegressor(max_depth=3, n_estimators=50))])"
$ws.Range("G2").Value = 0.122648122766744
$ws.Range("H2").Value = 0.991
